$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.465.25'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.919.10'
$ws.Range('E3').Value = '  +0.99%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  +0.67%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.09'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('E6').Value = '  +0.56%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4830'
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4078'
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08241'
$ws.Range('E9').Value = '  +2.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.27'
$ws.Range('E11').Value = '  -0.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.959.48'
$ws.Range('E12').Value = '  +5.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.069'
$ws.Range('E13').Value = '  +2.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.258'
$ws.Range('E14').Value = '  +2.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.59'
$ws.Range('E15').Value = '  +2.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06879'
$ws.Range('E16').Value = '  +2.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.010'
$ws.Range('E17').Value = '  +0.58%  '
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('E20').Value = '  +0.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '29.468.39'
$ws.Range('E21').Value = '  +0.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.667'
$ws.Range('E22').Value = '  +2.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.77'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.180'
$ws.Range('E24').Value = '  +1.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.143.41'
$ws.Range('E25').Value = '  +2.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.691'
$ws.Range('E26').Value = '  +10.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '156.22'
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.03'
$ws.Range('E28').Value = '  +1.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.118'
$ws.Range('E29').Value = '  +1.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.71'
$ws.Range('E30').Value = '  +1.89%  '
$ws.Range('E31').Value = '  -1.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09608'
$ws.Range('E32').Value = '  +1.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.664'
$ws.Range('E33').Value = '  +4.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.556'
$ws.Range('E34').Value = '  +0.39%  '
$ws.Range('E35').Value = '  -0.85%  '
$ws.Range('E36').Value = '  +1.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06112'
$ws.Range('E37').Value = '  +0.77%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.183'
$ws.Range('E38').Value = '  +0.72%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.082'
$ws.Range('E39').Value = '  +2.66%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5986'
$ws.Range('E40').Value = '  +2.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.84'
$ws.Range('E41').Value = '  +6.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1850'
$ws.Range('E42').Value = '  +0.42%  '
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.398'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.07601'
$ws.Range('E45').Value = '  -1.97%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.44'
$ws.Range('E46').Value = '  +1.66%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5604'
$ws.Range('E47').Value = '  +1.53%  '
$ws.Range('E48').Value = '  +1.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '118.17'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.430'
$ws.Range('E50').Value = '  +3.90%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.34'
$ws.Range('E51').Value = '  +0.38%  '
